$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$newValue = 85.77505782882612

$ws.Range("N2:N5").Value = $newValue
